# Swap the two "Periodo Mora" / "Valor Mora" data rows (E16/F16 <-> E17/F17)
# Before: E16=2412, F16=143600 | E17=2501, F17=14360
# After:  E16=2501, F16=14360  | E17=2412, F17=143600

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2501"
$ws.Range("F16").Value = 14360

$ws.Range("E17").Value = "2412"
$ws.Range("F17").Value = 143600
